$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "317.05"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "4.24%"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.87%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.165"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.38%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08250"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "5.09%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.148"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.62%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.031"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.36%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9277"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.92%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1024"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "5.06%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1888"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.60%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09408"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "9.34%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03588"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.14%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09926"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.26%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001435"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.41%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005653"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.82%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.467"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.02%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.143"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.04%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "11.38%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.51%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1331"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.53%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.180"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.00%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2193"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.37%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.08%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001247"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.07%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004733"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-6.55%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001252"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-21.90%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004506"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-5.22%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02002"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "8.97%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04943"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "4.40%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007936"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "4.76%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.14%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007847"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "1.55%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002114"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-4.21%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01174"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "4.06%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006510"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.88%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.04%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "34.19"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-27.40%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001903"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002103"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.04%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002003"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.04%"
